# Update the regression-results table ("CZG_PhenT_allScaled") with the
# re-run model numbers (PdeltaAIC added as a covariate for CG as well).
# Only the Estimate / SE / Chi2 / p-value columns move; every value is
# stored as TEXT (shared string), matching how the sheet already stores
# these look-like-numbers cells (no explicit cell style / number format).
#
# Plain `$ws.Range(...).Value = "0.442"` would get auto-coerced to a
# number by Excel, and forcing text via NumberFormat="@" (or a quote
# prefix) stamps a brand-new cell style onto the cell. Neither matches
# the source file, where these cells keep their original (default)
# style. Writing a `="text"` formula and then Copy/PasteSpecial-values
# collapses it back down to a literal shared-string value without ever
# touching the cell's style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Formula = '="' + $text + '"'
    $rng.Copy() | Out-Null
    $rng.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Row 2 - intrcpt
Set-TextValue 'B2' ' 0.013'
Set-TextValue 'D2' '0.442'
Set-TextValue 'F2' '0.6582'

# Row 3 - Diet_HCOherbivore
Set-TextValue 'B3' '-0.007'
Set-TextValue 'D3' '1.765'
Set-TextValue 'F3' '0.4138'

# Row 4 - Diet_HCOomnivore
Set-TextValue 'B4' ' 0.040'
Set-TextValue 'D4' '1.765'
Set-TextValue 'F4' '0.4138'

# Row 5 - Migratmigrant (Estimate/SE unchanged)
Set-TextValue 'D5' '0.151'
Set-TextValue 'F5' '0.6973'

# Row 6 - GenLength_y_IUCN.y
Set-TextValue 'B6' ' 0.001'
Set-TextValue 'C6' '0.002'
Set-TextValue 'D6' '0.323'
Set-TextValue 'F6' '0.5698'

# Row 7 - abs_lat (Estimate/SE unchanged)
Set-TextValue 'D7' '1.142'
Set-TextValue 'F7' '0.2851'

# Row 8 - Pvalue
Set-TextValue 'B8' '-0.044'
Set-TextValue 'D8' '1.100'
Set-TextValue 'F8' '0.2943'

$excel.CutCopyMode = $false
